{"js": "// Corre\u00e7\u00e3o no artefato 18\n// 1) \"Informar Resultado\" -> \"Contatar o candidato\" (the second \"Processo:\" occurrence,\n//    i.e. the one immediately followed by \"Evento: Candidato ... n\u00e3o solicita resultado\")\n// 2) \"Evento:\" run-group \"Candidato\" + \" n\u00e3o\" + \" solicita resultado\" -> single run\n//    \"Candidato n\u00e3o solicita resultado\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet processoPara = null;\nlet eventoPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"Processo: Informar Resultado\") {\n    const next = paragraphs.items[i + 1];\n    if (next) {\n      next.load(\"text\");\n      await context.sync();\n      if (next.text.indexOf(\"Evento:\") === 0 && next.text.indexOf(\"n\u00e3o solicita resultado\") !== -1) {\n        processoPara = para;\n        eventoPara = next;\n        break;\n      }\n    }\n  }\n}\n\nif (!processoPara || !eventoPara) {\n  throw new Error(\"Could not locate the target 'Informar Resultado' process block.\");\n}\n\n// --- Edit 1: replace \"Informar Resultado\" with \"Contatar o candidato\" ---\nconst procHits = processoPara.search(\"Informar Resultado\", { matchCase: true });\nprocHits.load(\"items\");\nawait context.sync();\nif (procHits.items.length === 0) {\n  throw new Error(\"Could not find 'Informar Resultado' text to replace.\");\n}\nprocHits.items[0].insertText(\"Contatar o candidato\", \"Replace\");\nawait context.sync();\n\n// --- Edit 2: merge the 3 runs into a single run with the full sentence ---\nconst eventHits = eventoPara.search(\"Candidato n\u00e3o solicita resultado\", { matchCase: true });\neventHits.load(\"items\");\nawait context.sync();\nif (eventHits.items.length === 0) {\n  throw new Error(\"Could not find 'Candidato n\u00e3o solicita resultado' text to normalize.\");\n}\neventHits.items[0].insertText(\"Candidato n\u00e3o solicita resultado\", \"Replace\");\nawait context.sync();\n", "ps1": "# Correcao no artefato 18\n# 1) \"Informar Resultado\" -> \"Contatar o candidato\" for the process block whose\n#    following \"Evento:\" paragraph is \"Candidato nao solicita resultado\"\n#    (there are two \"Informar Resultado\" processes in the doc; this is the\n#    second one).\n# 2) \"Evento:\" text \"Candidato\" + \" nao\" + \" solicita resultado\" (3 runs) is\n#    normalized into a single run \"Candidato nao solicita resultado\".\n\n$d = $word.ActiveDocument\n\n# --- Edit 1 -------------------------------------------------------------\n$rng = $d.Content\n$matchIndex = 0\n$target = $null\nwhile ($rng.Find.Execute(\"Informar Resultado\", $true)) {\n    $matchIndex = $matchIndex + 1\n    if ($matchIndex -eq 2) {\n        $target = $rng\n        break\n    }\n    $rng.Collapse(0)\n}\nif ($target -ne $null) {\n    $target.Text = \"Contatar o candidato\"\n}\n\n# --- Edit 2 -------------------------------------------------------------\n# Word's Find can match text spanning multiple runs; setting .Text on that\n# range collapses every run it covers into a single new run. Because the\n# combined text (\"Candidato\" + \" n\u00e3o\" + \" solicita resultado\") is already\n# identical to the desired final text, a direct same-value assignment is a\n# no-op, so we briefly round-trip through a placeholder to force the splice.\n$rng2 = $d.Content\nif ($rng2.Find.Execute(\"Candidato n\u00e3o solicita resultado\", $true)) {\n    $rng2.Text = \"Candidato n\u00e3o solicita resultado#__tmp__#\"\n    $rng3 = $d.Content\n    $rng3.Find.Execute(\"Candidato n\u00e3o solicita resultado#__tmp__#\", $true) | Out-Null\n    $rng3.Text = \"Candidato n\u00e3o solicita resultado\"\n}\n"}
